# Update the NATMI Slit3-Robo1 LR-pairs sheet with the re-run ("new TPM")
# values: every existing row (2-16) gets refreshed numbers, and the sending
# cluster x target cluster grid is extended from 3 target clusters to 4
# (adding "Resolving-Mac" as a target), producing 5 new rows (17-21) and
# growing the used range from A1:T16 to A1:T21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: row#, A..D (strings), E..T (numeric metrics)
$data = @(
    @(2, "ECs", "Slit3", "Robo1", "ECs", 3, 1, 2.384498666666666, 7.153495999999999, 0.03997130847982828, 0.03997130847982829, 3, 1, 0.7788713333333334, 2.336614, 0.02154486428780801, 0.02154486428780801, 1.857217655838222, 16.714958902544, 0.0008611764166040098, 0.0008611764166040098),
    @(3, "ECs", "Slit3", "Robo1", "FAPs", 3, 1, 2.384498666666666, 7.153495999999999, 0.03997130847982828, 0.03997130847982829, 3, 1, 29.811843, 89.435529, 0.8246446930529893, 0.8246446930529893, 71.08629988437598, 639.7766989593839, 0.03296212741227434, 0.03296212741227434),
    @(4, "ECs", "Slit3", "Robo1", "MuSCs", 3, 1, 2.384498666666666, 7.153495999999999, 0.03997130847982828, 0.03997130847982829, 3, 1, 5.508900000000001, 16.5267, 0.152385250035015, 0.152385250035015, 13.1359647048, 118.2236823432, 0.006091037836925349, 0.006091037836925351),
    @(5, "ECs", "Slit3", "Robo1", "Resolving-Mac", 3, 1, 2.384498666666666, 7.153495999999999, 0.03997130847982828, 0.03997130847982829, 1, 0.3333333333333333, 0.05152233333333334, 0.154567, 0.001425192624187658, 0.001425192624187658, 0.1228549351368889, 1.105694416232, [double]"5.696681402458086E-05", [double]"5.696681402458087E-05"),
    @(6, "FAPs", "Slit3", "Robo1", "ECs", 3, 1, 56.014214, 168.042642, 0.9389652669334476, 0.9389652669334477, 3, 1, 0.7788713333333334, 2.336614, 0.02154486428780801, 0.02154486428780801, 43.62786554379867, 392.650789894188, 0.02022987924704655, 0.02022987924704655),
    @(7, "FAPs", "Slit3", "Robo1", "FAPs", 3, 1, 56.014214, 168.042642, 0.9389652669334476, 0.9389652669334477, 3, 1, 29.811843, 89.435529, 0.8246446930529893, 0.8246446930529893, 1669.886953536402, 15028.98258182762, 0.774312724337751, 0.7743127243377511),
    @(8, "FAPs", "Slit3", "Robo1", "MuSCs", 3, 1, 56.014214, 168.042642, 0.9389652669334476, 0.9389652669334477, 3, 1, 5.508900000000001, 16.5267, 0.152385250035015, 0.152385250035015, 308.5767035046001, 2777.1903315414, 0.1430844569758481, 0.1430844569758481),
    @(9, "FAPs", "Slit3", "Robo1", "Resolving-Mac", 3, 1, 56.014214, 168.042642, 0.9389652669334476, 0.9389652669334477, 1, 0.3333333333333333, 0.05152233333333334, 0.154567, 0.001425192624187658, 0.001425192624187658, 2.885983005112667, 25.973847046014, 0.001338206372801945, 0.001338206372801945),
    @(10, "Inflammatory-Mac", "Slit3", "Robo1", "ECs", 2, 0.6666666666666666, 0.036048, 0.108144, 0.0006042719789376481, 0.0006042719789376483, 3, 1, 0.7788713333333334, 2.336614, 0.02154486428780801, 0.02154486428780801, 0.028076753824, 0.252690784416, [double]"1.301895777913681E-05", [double]"1.301895777913681E-05"),
    @(11, "Inflammatory-Mac", "Slit3", "Robo1", "FAPs", 2, 0.6666666666666666, 0.036048, 0.108144, 0.0006042719789376481, 0.0006042719789376483, 3, 1, 29.811843, 89.435529, 0.8246446930529893, 0.8246446930529893, 1.074657316464, 9.671915848175999, 0.0004983096805915593, 0.0004983096805915593),
    @(12, "Inflammatory-Mac", "Slit3", "Robo1", "MuSCs", 2, 0.6666666666666666, 0.036048, 0.108144, 0.0006042719789376481, 0.0006042719789376483, 3, 1, 5.508900000000001, 16.5267, 0.152385250035015, 0.152385250035015, 0.1985848272, 1.7872634448, [double]"9.208213659956685E-05", [double]"9.208213659956688E-05"),
    @(13, "Inflammatory-Mac", "Slit3", "Robo1", "Resolving-Mac", 2, 0.6666666666666666, 0.036048, 0.108144, 0.0006042719789376481, 0.0006042719789376483, 1, 0.3333333333333333, 0.05152233333333334, 0.154567, 0.001425192624187658, 0.001425192624187658, 0.001857277072, 0.016715493648, [double]"8.612039673852161E-07", [double]"8.612039673852163E-07"),
    @(14, "MuSCs", "Slit3", "Robo1", "ECs", 3, 1, 1.188166, 3.564498, 0.01991720539631685, 0.01991720539631685, 3, 1, 0.7788713333333334, 2.336614, 0.02154486428780801, 0.02154486428780801, 0.9254284366413335, 8.328855929772001, 0.0004291134872560438, 0.0004291134872560438),
    @(15, "MuSCs", "Slit3", "Robo1", "FAPs", 3, 1, 1.188166, 3.564498, 0.01991720539631685, 0.01991720539631685, 3, 1, 29.811843, 89.435529, 0.8246446930529893, 0.8246446930529893, 35.421418249938, 318.792764249442, 0.01642461773051905, 0.01642461773051905),
    @(16, "MuSCs", "Slit3", "Robo1", "MuSCs", 3, 1, 1.188166, 3.564498, 0.01991720539631685, 0.01991720539631685, 3, 1, 5.508900000000001, 16.5267, 0.152385250035015, 0.152385250035015, 6.545487677400001, 58.9093890966, 0.003035088324316494, 0.003035088324316494),
    @(17, "MuSCs", "Slit3", "Robo1", "Resolving-Mac", 3, 1, 1.188166, 3.564498, 0.01991720539631685, 0.01991720539631685, 1, 0.3333333333333333, 0.05152233333333334, 0.154567, 0.001425192624187658, 0.001425192624187658, 0.06121708470733334, 0.550953762366, [double]"2.83858542252614E-05", [double]"2.83858542252614E-05"),
    @(18, "Resolving-Mac", "Slit3", "Robo1", "ECs", 3, 1, 0.03233, 0.09698999999999999, 0.0005419472114695452, 0.0005419472114695452, 3, 1, 0.7788713333333334, 2.336614, 0.02154486428780801, 0.02154486428780801, 0.02518091020666667, 0.22662819186, [double]"1.167617912226734E-05", [double]"1.167617912226734E-05"),
    @(19, "Resolving-Mac", "Slit3", "Robo1", "FAPs", 3, 1, 0.03233, 0.09698999999999999, 0.0005419472114695452, 0.0005419472114695452, 3, 1, 29.811843, 89.435529, 0.8246446930529893, 0.8246446930529893, 0.9638168841899999, 8.67435195771, 0.0004469138918532266, 0.0004469138918532266),
    @(20, "Resolving-Mac", "Slit3", "Robo1", "MuSCs", 3, 1, 0.03233, 0.09698999999999999, 0.0005419472114695452, 0.0005419472114695452, 3, 1, 5.508900000000001, 16.5267, 0.152385250035015, 0.152385250035015, 0.178102737, 1.602924633, [double]"8.258476132556582E-05", [double]"8.258476132556582E-05"),
    @(21, "Resolving-Mac", "Slit3", "Robo1", "Resolving-Mac", 3, 1, 0.03233, 0.09698999999999999, 0.0005419472114695452, 0.0005419472114695452, 1, 0.3333333333333333, 0.05152233333333334, 0.154567, 0.001425192624187658, 0.001425192624187658, 0.001665717036666667, 0.01499145333, [double]"7.723791684854649E-07", [double]"7.723791684854649E-07")
)

foreach ($row in $data) {
    $r = [int]$row[0]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
